$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank spacer rows that separated each day's block (rows 9, 15,
# 21, 27 in the original layout). Deleting from the bottom up keeps the
# remaining row numbers stable while we work.
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(21).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(9).Delete()

# Restore the selection Excel leaves behind after this kind of edit.
$null = $ws.Range("H4:H28").Select()
